# [quickfort] add documentation and test data for config mode (#2070)
#
# The test fixture sheet lists one quickfort "mode" marker per row. This
# change inserts a new "#config hidden()" row for the new config mode,
# pushing the #meta/#notes rows down and appending a new #notes row at
# the bottom so every marker keeps its original relative ordering:
#
#   #dig, #build, #place, #zone, #query, #config (NEW), #ignore,
#   #aliases, #meta, #notes
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the tail of the list down one row to make room for the new
# "#config hidden()" entry, then fill in the new row.
$ws.Range("A6").Value = "#config hidden()"
$ws.Range("A7").Value = "#ignore"
$ws.Range("A8").Value = "#aliases"
$ws.Range("A9").Value = "#meta hidden()"
$ws.Range("A10").Value = "#notes hidden()"

# New row 10 should pick up the same cell formatting as the rest of the
# column (style index 1 in the original workbook) rather than defaulting
# to unformatted.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
